# Rename the two "Include" worksheets
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(2).Name = "Include #0"
$wb.Worksheets.Item(3).Name = "Include #1"

# Update values on the Metadata sheet
$ws = $wb.Worksheets.Item(1)

# URL (row 2) - pythia -> cicada
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/VaccineCodesCvxMvx"

# Date (row 8) - refreshed publication date
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# Insert a "Jurisdiction" property row right after "Contact" (row 10), which
# pushes Description/Purpose/Copyright/Immutable down by one row each.
# Copy precise A:B ranges (not whole rows) from bottom to top so the
# existing cell styles are preserved instead of generating new style ids.
$ws.Range("B15").ClearContents()
$ws.Range("A14:B14").Copy($ws.Range("A15:B15"))

$ws.Range("B14").ClearContents()
$ws.Range("A13:B13").Copy($ws.Range("A14:B14"))

$ws.Range("B13").ClearContents()
$ws.Range("A12:B12").Copy($ws.Range("A13:B13"))

$ws.Range("B12").ClearContents()
$ws.Range("A11:B11").Copy($ws.Range("A12:B12"))

# Write the new "Jurisdiction" row
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

$wb.Save()
